# Insert a new weekly price record for "Espinaca" (Vega Central Mapocho de
# Santiago) above the existing row 405, pushing all subsequent rows down by
# one (old row 405 -> new row 406, ..., old row 468 -> new row 469).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 405; everything from 405 downward shifts
# down to make room (old row 405 becomes row 406, etc.).
$ws.Rows.Item(405).Insert()

# Populate the newly inserted row 405 with the new record's data.
$ws.Cells.Item(405, 1).Value = 9
$ws.Cells.Item(405, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(405, 3).Value = "Metropolitana"
$ws.Cells.Item(405, 4).Value = 44776
$ws.Cells.Item(405, 5).Value = 13
$ws.Cells.Item(405, 6).Value = 100112012
$ws.Cells.Item(405, 7).Value = "Espinaca"
$ws.Cells.Item(405, 8).Value = "Sin especificar"
$ws.Cells.Item(405, 9).Value = "Primera"
$ws.Cells.Item(405, 10).Value = 160
$ws.Cells.Item(405, 11).Value = 8000
$ws.Cells.Item(405, 12).Value = 9000
$ws.Cells.Item(405, 13).Value = 8500
$ws.Cells.Item(405, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(405, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(405, 16).Value = 850
$ws.Cells.Item(405, 17).Value = 10
$ws.Cells.Item(405, 18).Value = "Hortaliza"
